$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

$ws.Range("B12").Value = 64
$ws.Range("C12").Value = -14
$ws.Range("E12").Value = "50 / 112"
